$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value for every data row (rows 2-496).
# All of them are being bumped from 45178 (2023-09-09) to 45179 (2023-09-10).
$range = $ws.Range("C2:C496")
$range.Value = 45179
